# Natmi following Dr Hou advice
# Update the Fgf10-Fgfrl1 LR-pair computed values: the number of
# ligand-/receptor-expressing cells changed from 1 to 3, which cascades
# into the total/average expression and the derived-specificity columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.528376666666667
$ws.Range("H2").Value = 4.58513
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7894570000000001
$ws.Range("N2").Value = 2.368371
$ws.Range("O2").Value = 0.09838606084581891
$ws.Range("P2").Value = 0.09838606084581894
$ws.Range("Q2").Value = 1.206587658136667
$ws.Range("R2").Value = 10.85928892323
$ws.Range("S2").Value = 0.09838606084581891
$ws.Range("T2").Value = 0.09838606084581894

# --- Row 3 ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.528376666666667
$ws.Range("H3").Value = 4.58513
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.625751333333334
$ws.Range("N3").Value = 16.877254
$ws.Range("O3").Value = 0.7011091332204036
$ws.Range("P3").Value = 0.7011091332204038
$ws.Range("Q3").Value = 8.598267070335556
$ws.Range("R3").Value = 77.38440363302001
$ws.Range("S3").Value = 0.7011091332204036
$ws.Range("T3").Value = 0.7011091332204038

# --- Row 4 ---
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.528376666666667
$ws.Range("H4").Value = 4.58513
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.608865333333333
$ws.Range("N4").Value = 4.826596
$ws.Range("O4").Value = 0.2005048059337773
$ws.Range("P4").Value = 0.2005048059337774
$ws.Range("Q4").Value = 2.458952235275556
$ws.Range("R4").Value = 22.13057011748
$ws.Range("S4").Value = 0.2005048059337773
$ws.Range("T4").Value = 0.2005048059337774
